$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("text_coercion")

# Add a new row of test data: a numeric "student number" value paired with
# a text label explaining what is being tested.
$ws.Range("A9").Value = 36436153
$ws.Range("B9").Value = "student number"

# Move/update the active selection to A9 (it previously pointed at B9).
$ws.Range("A9").Select()
